$wb = $excel.ActiveWorkbook

# OFF sheet - row 3 (R)
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 253
$wsOff.Range("C3").Value = 157
$wsOff.Range("D3").Value = 52
$wsOff.Range("F3").Value = 6
$wsOff.Range("G3").Value = 4

# DEF sheet - row 3 (R)
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 253
$wsDef.Range("C3").Value = 188
$wsDef.Range("D3").Value = 46
$wsDef.Range("E3").Value = 20
